$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row2 = New-Object 'object[,]' 1,10
$row2[0,0] = -0.5351204465965399
$row2[0,1] = -1.66950562919271
$row2[0,2] = -0.2895456687149427
$row2[0,3] = -0.6733414736251095
$row2[0,4] = 0.0115444171491989
$row2[0,5] = -0.6804809672324722
$row2[0,6] = -0.4410326232298434
$row2[0,7] = -0.4043046464706727
$row2[0,8] = 0.4221894279166428
$row2[0,9] = -0.2336087822049224
$ws.Range("B2:K2").Value = $row2

$row3 = New-Object 'object[,]' 1,10
$row3[0,0] = 0.01136921538350649
$row3[0,1] = -0.3724265895266604
$row3[0,2] = 0.3124593012476481
$row3[0,3] = -0.3795660831340231
$row3[0,4] = -0.1401177391313942
$row3[0,5] = -0.1033897623722235
$row3[0,6] = 0.723104312015092
$row3[0,7] = 0.06730610189352677
$row3[0,8] = -0.5908441378320941
$row3[0,9] = -0.191729189980311
$ws.Range("B3:K3").Value = $row3

$row4 = New-Object 'object[,]' 1,10
$row4[0,0] = 0.9526635850093691
$row4[0,1] = 0.2606382006276979
$row4[0,2] = 0.5000865446303268
$row4[0,3] = 0.5368145213894975
$row4[0,4] = 1.363308595776813
$row4[0,5] = 0.7075103856552477
$row4[0,6] = 0.0493601459296269
$row4[0,7] = 0.44847509378141
$row4[0,8] = 0.3870385516598326
$row4[0,9] = -0.04712738345727097
$ws.Range("B4:K4").Value = $row4

$row5 = New-Object 'object[,]' 1,10
$row5[0,0] = 0.4274989419678774
$row5[0,1] = 0.4642269187270481
$row5[0,2] = 1.290720993114364
$row5[0,3] = 0.6349227829927984
$row5[0,4] = -0.02322745673282245
$row5[0,5] = 0.3758874911189606
$row5[0,6] = 0.3144509489973832
$row5[0,7] = -0.1197149861197203
$row5[0,8] = 0.8590864059593566
$row5[0,9] = 0.6319705015114304
$ws.Range("B5:K5").Value = $row5

$row6 = New-Object 'object[,]' 1,10
$row6[0,0] = 1.235029948750828
$row6[0,1] = 0.5792317386292632
$row6[0,2] = -0.07891850109635756
$row6[0,3] = 0.3201964467554255
$row6[0,4] = 0.2587599046338481
$row6[0,5] = -0.1754060304832554
$row6[0,6] = 0.8033953615958215
$row6[0,7] = 0.5762794571478953
$row6[0,8] = -0.05718027227819322
$row6[0,9] = 0.3702166863774111
$ws.Range("B6:K6").Value = $row6

$row7 = New-Object 'object[,]' 1,10
$row7[0,0] = -0.4329776612703231
$row7[0,1] = -0.03386271341853997
$row7[0,2] = -0.09529925554011737
$row7[0,3] = -0.529465190657221
$row7[0,4] = 0.449336201421856
$row7[0,5] = 0.2222202969739298
$row7[0,6] = -0.4112394324521587
$row7[0,7] = 0.01615752620344563
$row7[0,8] = -0.1501794284847013
$row7[0,9] = -0.02746841204387546
$ws.Range("B7:K7").Value = $row7

$row8 = New-Object 'object[,]' 1,10
$row8[0,0] = -0.1856174723396913
$row8[0,1] = -0.6197834074567948
$row8[0,2] = 0.3590179846222821
$row8[0,3] = 0.1319020801743559
$row8[0,4] = -0.5015576492517326
$row8[0,5] = -0.07416069059612829
$row8[0,6] = -0.2404976452842752
$row8[0,7] = -0.1177866288434494
$row8[0,8] = -0.07055289228830908
$row8[0,9] = -0.4671716238107607
$ws.Range("B8:K8").Value = $row8

$row9 = New-Object 'object[,]' 1,10
$row9[0,0] = 0.4273407187267424
$row9[0,1] = 0.2002248142788162
$row9[0,2] = -0.4332349151472724
$row9[0,3] = -0.005837956491668017
$row9[0,4] = -0.1721749111798149
$row9[0,5] = -0.0494638947389891
$row9[0,6] = -0.002230158183848807
$row9[0,7] = -0.3988488897063004
$row9[0,8] = -0.1447968545825803
$row9[0,9] = 0.07532753529099229
$ws.Range("B9:K9").Value = $row9

$row10 = New-Object 'object[,]' 1,10
$row10[0,0] = -0.5318964931771777
$row10[0,1] = -0.1044995345215733
$row10[0,2] = -0.2708364892097202
$row10[0,3] = -0.1481254727688944
$row10[0,4] = -0.1008917362137541
$row10[0,5] = -0.4975104677362057
$row10[0,6] = -0.2434584326124856
$row10[0,7] = -0.02333404273891299
$row10[0,8] = -0.2591224913255812
$row10[0,9] = -0.4380952487963659
$ws.Range("B10:K10").Value = $row10

$row11 = New-Object 'object[,]' 1,10
$row11[0,0] = -0.1938269109680474
$row11[0,1] = -0.07111589452722158
$row11[0,2] = -0.02388215797208129
$row11[0,3] = -0.4205008894945329
$row11[0,4] = -0.1664488543708128
$row11[0,5] = 0.0536755355027598
$row11[0,6] = -0.1821129130839084
$row11[0,7] = -0.3610856705546931
$row11[0,8] = -0.3716462008140141
$row11[0,9] = -0.5793653109721442
$ws.Range("B11:K11").Value = $row11

$row12 = New-Object 'object[,]' 1,10
$row12[0,0] = 0.1476338940440795
$row12[0,1] = -0.2489848374783721
$row12[0,2] = 0.005067197645347965
$row12[0,3] = 0.2251915875189206
$row12[0,4] = -0.0105968610677476
$row12[0,5] = -0.1895696185385323
$row12[0,6] = -0.2001301487978533
$row12[0,7] = -0.4078492589559834
$row12[0,8] = -0.1906403594810787
$row12[0,9] = 0.1190458097769828
$ws.Range("B12:K12").Value = $row12

$row13 = New-Object 'object[,]' 1,10
$row13[0,0] = -0.0323979044984018
$row13[0,1] = 0.1877264853751708
$row13[0,2] = -0.04806196321149736
$row13[0,3] = -0.2270347206822821
$row13[0,4] = -0.2375952509416031
$row13[0,5] = -0.4453143610997332
$row13[0,6] = -0.2281054616248284
$row13[0,7] = 0.08158070763323305
$row13[0,8] = -0.07272342619877098
$row13[0,9] = 0.5777029950204122
$ws.Range("B13:K13").Value = $row13

$row14 = New-Object 'object[,]' 1,10
$row14[0,0] = -0.06996447561954
$row14[0,1] = -0.2489372330903247
$row14[0,2] = -0.2594977633496457
$row14[0,3] = -0.4672168735077758
$row14[0,4] = -0.2500079740328711
$row14[0,5] = 0.05967819522519041
$row14[0,6] = -0.09462593860681362
$row14[0,7] = 0.5558004826123696
$row14[0,8] = 0.3396354339941604
$row14[0,9] = -0.0673936950407959
$ws.Range("B14:K14").Value = $row14

$row15 = New-Object 'object[,]' 1,10
$row15[0,0] = -0.1713918715036764
$row15[0,1] = -0.3791109816618064
$row15[0,2] = -0.1619020821869017
$row15[0,3] = 0.1477840870711598
$row15[0,4] = -0.006520046760844223
$row15[0,5] = 0.643906374458339
$row15[0,6] = 0.4277413258401298
$row15[0,7] = 0.0207121968051735
$row15[0,8] = 0.5197544139825933
$row15[0,9] = 0.3609055008270807
$ws.Range("B15:K15").Value = $row15

$row16 = New-Object 'object[,]' 1,10
$row16[0,0] = -0.06818896562035748
$row16[0,1] = 0.241497203637704
$row16[0,2] = 0.08719306980569996
$row16[0,3] = 0.7376194910248832
$row16[0,4] = 0.521454442406674
$row16[0,5] = 0.1144253133717177
$row16[0,6] = 0.6134675305491375
$row16[0,7] = 0.4546186173936249
$row16[0,8] = 0.5107824383638689
$row16[0,9] = 2.760585277975261
$ws.Range("B16:K16").Value = $row16

$row17 = New-Object 'object[,]' 1,10
$row17[0,0] = 0.2497007499081394
$row17[0,1] = 0.09539661607613537
$row17[0,2] = 0.7458230372953185
$row17[0,3] = 0.5296579886771094
$row17[0,4] = 0.1226288596421531
$row17[0,5] = 0.6216710768195729
$row17[0,6] = 0.4628221636640603
$row17[0,7] = 0.5189859846343043
$row17[0,8] = 2.768788824245696
$row17[0,9] = 10.23793915510299
$ws.Range("B17:K17").Value = $row17

$row18 = New-Object 'object[,]' 1,10
$row18[0,0] = 0.09280705542466716
$row18[0,1] = 0.7432334766438504
$row18[0,2] = 0.5270684280256412
$row18[0,3] = 0.1200392989906849
$row18[0,4] = 0.6190815161681047
$row18[0,5] = 0.4602326030125921
$row18[0,6] = 0.5163964239828361
$row18[0,7] = 2.766199263594229
$row18[0,8] = 10.23534959445152
$row18[0,9] = -7.935912205685947
$ws.Range("B18:K18").Value = $row18

$row19 = New-Object 'object[,]' 1,10
$row19[0,0] = 0.7497668092269023
$row19[0,1] = 0.5336017606086931
$row19[0,2] = 0.1265726315737368
$row19[0,3] = 0.6256148487511566
$row19[0,4] = 0.466765935595644
$row19[0,5] = 0.522929756565888
$row19[0,6] = 2.77273259617728
$row19[0,7] = 10.24188292703457
$row19[0,8] = -7.929378873102896
$row19[0,9] = 0.2102926738762539
$ws.Range("B19:K19").Value = $row19

$row20 = New-Object 'object[,]' 1,10
$row20[0,0] = 0.4223850656296224
$row20[0,1] = 0.01535593659466611
$row20[0,2] = 0.5143981537720859
$row20[0,3] = 0.3555492406165733
$row20[0,4] = 0.4117130615868174
$row20[0,5] = 2.661515901198209
$row20[0,6] = 10.1306662320555
$row20[0,7] = -8.040595568081965
$row20[0,8] = 0.09907597889718328
$row20[0,9] = 2.246646450696576
$ws.Range("B20:K20").Value = $row20

$row21 = New-Object 'object[,]' 1,10
$row21[0,0] = -0.02639020739223796
$row21[0,1] = 0.4726520097851818
$row21[0,2] = 0.3138030966296693
$row21[0,3] = 0.3699669175999133
$row21[0,4] = 2.619769757211305
$row21[0,5] = 10.0889200880686
$row21[0,6] = -8.082341712068869
$row21[0,7] = 0.05732983491027921
$row21[0,8] = 2.204900306709672
$row21[0,9] = -1.235129679813658
$ws.Range("B21:K21").Value = $row21

$row22 = New-Object 'object[,]' 1,10
$row22[0,0] = 0.4979670725178967
$row22[0,1] = 0.3391181593623842
$row22[0,2] = 0.3952819803326282
$row22[0,3] = 2.64508481994402
$row22[0,4] = 10.11423515080131
$row22[0,5] = -8.057026649336155
$row22[0,6] = 0.0826448976429941
$row22[0,7] = 2.230215369442386
$row22[0,8] = -1.209814617080943
$row22[0,9] = -1.270988795495144
$ws.Range("B22:K22").Value = $row22

$row23 = New-Object 'object[,]' 1,10
$row23[0,0] = 0.343156824405298
$row23[0,1] = 0.3993206453755421
$row23[0,2] = 2.649123484986935
$row23[0,3] = 10.11827381584423
$row23[0,4] = -8.052987984293241
$row23[0,5] = 0.08668356268590799
$row23[0,6] = 2.2342540344853
$row23[0,7] = -1.20577595203803
$row23[0,8] = -1.266950130452231
$row23[0,9] = 0.7745058067040239
$ws.Range("B23:K23").Value = $row23

$row24 = New-Object 'object[,]' 1,10
$row24[0,0] = 0.2804435086845197
$row24[0,1] = 2.530246348295912
$row24[0,2] = 9.999396679153206
$row24[0,3] = -8.171865120984263
$row24[0,4] = -0.03219357400511441
$row24[0,5] = 2.115376897794278
$row24[0,6] = -1.324653088729052
$row24[0,7] = -1.385827267143253
$row24[0,8] = 0.6556286700130015
$row24[0,9] = 0.07303413297936051
$ws.Range("B24:K24").Value = $row24

$row25 = New-Object 'object[,]' 1,10
$row25[0,0] = 2.48932270964054
$row25[0,1] = 9.958473040497832
$row25[0,2] = -8.212788759639636
$row25[0,3] = -0.07311721266048643
$row25[0,4] = 2.074453259138906
$row25[0,5] = -1.365576727384424
$row25[0,6] = -1.426750905798625
$row25[0,7] = 0.6147050313576294
$row25[0,8] = 0.03211049432398849
$row25[0,9] = 0.1277855351333463
$ws.Range("B25:K25").Value = $row25

$row26 = New-Object 'object[,]' 1,10
$row26[0,0] = 9.643547872076862
$row26[0,1] = -8.527713928060606
$row26[0,2] = -0.388042381081458
$row26[0,3] = 1.759528090717934
$row26[0,4] = -1.680501895805395
$row26[0,5] = -1.741676074219596
$row26[0,6] = 0.2997798629366579
$row26[0,7] = -0.2828146740969831
$row26[0,8] = -0.1871396332876253
$row26[0,9] = -0.4475586702863481
$ws.Range("B26:K26").Value = $row26

$row27 = New-Object 'object[,]' 1,10
$row27[0,0] = -9.584088888243137
$row27[0,1] = -1.444417341263988
$row27[0,2] = 0.7031531305354048
$row27[0,3] = -2.736876855987925
$row27[0,4] = -2.798051034402126
$row27[0,5] = -0.7565950972458717
$row27[0,6] = -1.339189634279513
$row27[0,7] = -1.243514593470155
$row27[0,8] = -1.503933630468878
$row27[0,9] = -0.6684786343103865
$ws.Range("B27:K27").Value = $row27

$row28 = New-Object 'object[,]' 1,10
$row28[0,0] = -0.573770965293057
$row28[0,1] = 1.573799506506335
$row28[0,2] = -1.866230480016994
$row28[0,3] = -1.927404658431195
$row28[0,4] = 0.1140512787250589
$row28[0,5] = -0.4685432583085821
$row28[0,6] = -0.3728682174992243
$row28[0,7] = -0.6332872544979471
$row28[0,8] = 0.2021677416605441
$row28[0,9] = 0.1020883817579226
$ws.Range("B28:K28").Value = $row28

$row29 = New-Object 'object[,]' 1,10
$row29[0,0] = 1.563148290176452
$row29[0,1] = -1.876881696346878
$row29[0,2] = -1.938055874761079
$row29[0,3] = 0.1034000623951754
$row29[0,4] = -0.4791944746384656
$row29[0,5] = -0.3835194338291078
$row29[0,6] = -0.6439384708278306
$row29[0,7] = 0.1915165253306606
$row29[0,8] = 0.09143716542803909
$row29[0,9] = -0.2623493406516572
$ws.Range("B29:K29").Value = $row29

$row30 = New-Object 'object[,]' 1,10
$row30[0,0] = -1.985496228563019
$row30[0,1] = -2.04667040697722
$row30[0,2] = -0.005214469820965406
$row30[0,3] = -0.5878090068546065
$row30[0,4] = -0.4921339660452486
$row30[0,5] = -0.7525530030439714
$row30[0,6] = 0.08290199311451979
$row30[0,7] = -0.01717736678810172
$row30[0,8] = -0.370963872867798
$row30[0,9] = -0.2227864824353526
$ws.Range("B30:K30").Value = $row30

$row31 = New-Object 'object[,]' 1,10
$row31[0,0] = -1.838568686009481
$row31[0,1] = 0.2028872511467736
$row31[0,2] = -0.3797072858868674
$row31[0,3] = -0.2840322450775096
$row31[0,4] = -0.5444512820762324
$row31[0,5] = 0.2910037140822588
$row31[0,6] = 0.1909243541796373
$row31[0,7] = -0.1628621519000589
$row31[0,8] = -0.01468476146761361
$row31[0,9] = 0.1427612275365414
$ws.Range("B31:K31").Value = $row31

$row32 = New-Object 'object[,]' 1,10
$row32[0,0] = 0.319385421520574
$row32[0,1] = -0.263209115513067
$row32[0,2] = -0.1675340747037092
$row32[0,3] = -0.427953111702432
$row32[0,4] = 0.4075018844560592
$row32[0,5] = 0.3074225245534377
$row32[0,6] = -0.04636398152625851
$row32[0,7] = 0.1018134089061868
$row32[0,8] = 0.2592593979103418
$row32[0,9] = -0.241400058615729
$ws.Range("B32:K32").Value = $row32

$row33 = New-Object 'object[,]' 1,10
$row33[0,0] = -0.2177157015159319
$row33[0,1] = -0.1220406607065741
$row33[0,2] = -0.3824596977052969
$row33[0,3] = 0.4529952984531944
$row33[0,4] = 0.3529159385505728
$row33[0,5] = -0.0008705675291234075
$row33[0,6] = 0.1473068229033219
$row33[0,7] = 0.3047528119074769
$row33[0,8] = -0.1959066446185939
$row33[0,9] = 0.08445123591687528
$ws.Range("B33:K33").Value = $row33

$row34 = New-Object 'object[,]' 1,10
$row34[0,0] = -0.1395947820665385
$row34[0,1] = -0.4000138190652613
$row34[0,2] = 0.4354411770932299
$row34[0,3] = 0.3353618171906084
$row34[0,4] = -0.01842468888908786
$row34[0,5] = 0.1297527015433575
$row34[0,6] = 0.2871986905475125
$row34[0,7] = -0.2134607659785583
$row34[0,8] = 0.06689711455691082
$row34[0,9] = -0.1050777432881008
$ws.Range("B34:K34").Value = $row34

$row35 = New-Object 'object[,]' 1,10
$row35[0,0] = -0.3119065001142551
$row35[0,1] = 0.5235484960442361
$row35[0,2] = 0.4234691361416146
$row35[0,3] = 0.06968263006191837
$row35[0,4] = 0.2178600204943637
$row35[0,5] = 0.3753060094985187
$row35[0,6] = -0.1253534470275521
$row35[0,7] = 0.155004433507917
$row35[0,8] = -0.01697042433709459
$row35[0,9] = 0.2888921154092369
$ws.Range("B35:K35").Value = $row35

$row36 = New-Object 'object[,]' 1,10
$row36[0,0] = 0.7021231295320197
$row36[0,1] = 0.6020437696293982
$row36[0,2] = 0.248257263549702
$row36[0,3] = 0.3964346539821473
$row36[0,4] = 0.5538806429863024
$row36[0,5] = 0.0532211864602315
$row36[0,6] = 0.3335790669957007
$row36[0,7] = 0.161604209150689
$row36[0,8] = 0.4674667488970205
$row36[0,9] = -0.1448632037902657
$ws.Range("B36:K36").Value = $row36

$row37 = New-Object 'object[,]' 1,10
$row37[0,0] = 1.514070997382048
$row37[0,1] = 1.160284491302352
$row37[0,2] = 1.308461881734797
$row37[0,3] = 1.465907870738952
$row37[0,4] = 0.9652484142128814
$row37[0,5] = 1.245606294748351
$row37[0,6] = 1.073631436903339
$row37[0,7] = 1.379493976649671
$row37[0,8] = 0.7671640239623843
$row37[0,9] = 1.455535409161496
$ws.Range("B37:K37").Value = $row37

$row38 = New-Object 'object[,]' 1,10
$row38[0,0] = 0.2163102553365951
$row38[0,1] = 0.3644876457690405
$row38[0,2] = 0.5219336347731955
$row38[0,3] = 0.02127417824712469
$row38[0,4] = 0.3016320587825939
$row38[0,5] = 0.1296572009375822
$row38[0,6] = 0.4355197406839137
$row38[0,7] = -0.1768102120033725
$row38[0,8] = 0.511561173195739
$row38[0,9] = 0.2348700177716323
$ws.Range("B38:K38").Value = $row38

$row39 = New-Object 'object[,]' 1,9
$row39[0,0] = 0.3684555432821496
$row39[0,1] = 0.5259015322863045
$row39[0,2] = 0.0252420757602338
$row39[0,3] = 0.305599956295703
$row39[0,4] = 0.1336250984506913
$row39[0,5] = 0.4394876381970228
$row39[0,6] = -0.1728423144902634
$row39[0,7] = 0.5155290707088481
$row39[0,8] = 0.2388379152847414
$ws.Range("B39:J39").Value = $row39

$row40 = New-Object 'object[,]' 1,8
$row40[0,0] = 0.661541622456546
$row40[0,1] = 0.1608821659304752
$row40[0,2] = 0.4412400464659443
$row40[0,3] = 0.2692651886209327
$row40[0,4] = 0.5751277283672642
$row40[0,5] = -0.03720222432002201
$row40[0,6] = 0.6511691608790895
$row40[0,7] = 0.3744780054549828
$ws.Range("B40:I40").Value = $row40

$row41 = New-Object 'object[,]' 1,7
$row41[0,0] = -0.07992401592518952
$row41[0,1] = 0.2004338646102796
$row41[0,2] = 0.028459006765268
$row41[0,3] = 0.3343215465115995
$row41[0,4] = -0.2780084061756867
$row41[0,5] = 0.4103629790234248
$row41[0,6] = 0.1336718235993181
$ws.Range("B41:H41").Value = $row41

$row42 = New-Object 'object[,]' 1,6
$row42[0,0] = 0.1551026493581833
$row42[0,1] = -0.01687220848682837
$row42[0,2] = 0.2889903312595031
$row42[0,3] = -0.3233396214277831
$row42[0,4] = 0.3650317637713285
$row42[0,5] = 0.08834060834722172
$ws.Range("B42:G42").Value = $row42

$row43 = New-Object 'object[,]' 1,5
$row43[0,0] = -0.08373363042288225
$row43[0,1] = 0.2221289093234493
$row43[0,2] = -0.3902010433638369
$row43[0,3] = 0.2981703418352746
$row43[0,4] = 0.02147918641116785
$ws.Range("B43:F43").Value = $row43

$row44 = New-Object 'object[,]' 1,4
$row44[0,0] = 0.1925427069667326
$row44[0,1] = -0.4197872457205535
$row44[0,2] = 0.268584139478558
$row44[0,3] = -0.00810701594554874
$ws.Range("B44:E44").Value = $row44

$row45 = New-Object 'object[,]' 1,3
$row45[0,0] = -0.4379379024501944
$row45[0,1] = 0.2504334827489171
$row45[0,2] = -0.02625767267518964
$ws.Range("B45:D45").Value = $row45

$row46 = New-Object 'object[,]' 1,2
$row46[0,0] = 0.2324016585002178
$row46[0,1] = -0.04428949692388896
$ws.Range("B46:C46").Value = $row46

$row47 = New-Object 'object[,]' 1,1
$row47[0,0] = -0.09587373626955231
$ws.Range("B47:B47").Value = $row47
